$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 18: update D18
$ws.Range("D18").Value = 0.6447426901493167

# Row 19: update C19, add D19
$ws.Range("C19").Value = 0.2386249091493167
$ws.Range("D19").Value = 0.597740902

# Row 20: update B20, add C20
$ws.Range("B20").Value = -0.0107480648506833
$ws.Range("C20").Value = 0.042359665
